$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# Update C13: new "Azami" (max) value for the inbound SWIFT row
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"

# Clear the benchmark values that are no longer available for SENET rows
$ws.Range("K24").ClearContents()
$ws.Range("K25").ClearContents()
